# Auto-generated Excel COM-interop script
# Applies updated pricing/profit values to the Ragnarok_Profits workbook
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR sheets,
# matching the scheduled runner's refreshed market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1741.5883
$ws.Range("I4").Value = 1538.75
$ws.Range("K4").Value = 1538.75
$ws.Range("M4").Value = -1424.75
$ws.Range("H55").Value = 1071.4615
$ws.Range("I55").Value = 214
$ws.Range("K55").Value = 214
$ws.Range("M55").Value = 0
$ws.Range("H132").Value = 2318.9714
$ws.Range("I132").Value = 2328.5
$ws.Range("K132").Value = 6985.5
$ws.Range("M132").Value = -4455.5
$ws.Range("H138").Value = 4830.8184
$ws.Range("J138").Value = 5971.8438
$ws.Range("L138").Value = 17915.5314
$ws.Range("N138").Value = -28195.5314

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12205.947
$ws.Range("I32").Value = 11769.56
$ws.Range("K32").Value = 11769.56
$ws.Range("M32").Value = -11482.56
$ws.Range("H60").Value = 89499.75
$ws.Range("I60").Value = 89499.75
$ws.Range("K60").Value = 89499.75
$ws.Range("M60").Value = -88766.75
$ws.Range("H110").Value = 4589.375
$ws.Range("I110").Value = 4248.6665
$ws.Range("J110").Value = 5611.5
$ws.Range("K110").Value = 4248.6665
$ws.Range("L110").Value = 5611.5
$ws.Range("M110").Value = -2203.6665
$ws.Range("N110").Value = -9701.5
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 65084.547
$ws.Range("I86").Value = 113576.5
$ws.Range("J86").Value = 6894.2
$ws.Range("K86").Value = 113576.5
$ws.Range("L86").Value = 6894.2
$ws.Range("M86").Value = -112453.5
$ws.Range("N86").Value = -9140.200000000001
$ws.Range("H89").Value = 65084.547
$ws.Range("I89").Value = 113576.5
$ws.Range("J89").Value = 6894.2
$ws.Range("K89").Value = 567882.5
$ws.Range("L89").Value = 34471
$ws.Range("M89").Value = -562266.5
$ws.Range("N89").Value = -45703

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 11112330
$ws.Range("I16").Value = 14286496
$ws.Range("J16").Value = 2750
$ws.Range("K16").Value = 14286496
$ws.Range("L16").Value = 2750
$ws.Range("M16").Value = -14286209
$ws.Range("N16").Value = -3324
$ws.Range("H58").Value = 1922.8823
$ws.Range("I58").Value = 1271.8182
$ws.Range("K58").Value = 1271.8182
$ws.Range("M58").Value = -1068.8182
$ws.Range("H94").Value = 1588.75
$ws.Range("I94").Value = 1395.5
$ws.Range("J94").Value = 1782
$ws.Range("K94").Value = 1395.5
$ws.Range("L94").Value = 1782
$ws.Range("M94").Value = -944.5
$ws.Range("N94").Value = -2684
$ws.Range("H113").Value = 11112330
$ws.Range("I113").Value = 14286496
$ws.Range("J113").Value = 2750
$ws.Range("K113").Value = 14286496
$ws.Range("L113").Value = 2750
$ws.Range("M113").Value = -14284326
$ws.Range("N113").Value = -7090
$ws.Range("H122").Value = 3422.9546
$ws.Range("I122").Value = 1718.5714
$ws.Range("J122").Value = 6405.625
$ws.Range("K122").Value = 5155.7142
$ws.Range("L122").Value = 19216.875
$ws.Range("M122").Value = -2705.7142
$ws.Range("N122").Value = -24116.875
$ws.Range("H134").Value = 2008.2858
$ws.Range("I134").Value = 2038.2222
$ws.Range("K134").Value = 6114.6666
$ws.Range("M134").Value = -3579.6666
$ws.Range("H136").Value = 1922.8823
$ws.Range("I136").Value = 1271.8182
$ws.Range("K136").Value = 3815.4546
$ws.Range("M136").Value = -1265.4546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 4014.8823
$ws.Range("I132").Value = 3855.2727
$ws.Range("J132").Value = 4091.2173
$ws.Range("K132").Value = 34697.4543
$ws.Range("L132").Value = 36820.9557
$ws.Range("M132").Value = -32167.4543
$ws.Range("N132").Value = -41880.9557
$ws.Range("H137").Value = 6785.5107
$ws.Range("I137").Value = 3318.85
$ws.Range("J137").Value = 9353.406999999999
$ws.Range("K137").Value = 9956.549999999999
$ws.Range("L137").Value = 28060.221
$ws.Range("M137").Value = -4856.549999999999
$ws.Range("N137").Value = -38260.221
$ws.Range("H139").Value = 4321.0356
$ws.Range("I139").Value = 2339.9443
$ws.Range("J139").Value = 7887
$ws.Range("K139").Value = 7019.8329
$ws.Range("L139").Value = 23661
$ws.Range("M139").Value = -1879.8329
$ws.Range("N139").Value = -33941

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 5002500
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 5002500
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 5002500
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -5002836
$ws.Range("H92").Value = 20749.6
$ws.Range("J92").Value = 20749.6
$ws.Range("L92").Value = 20749.6
$ws.Range("N92").Value = -24493.6
$ws.Range("H97").Value = 1027.7778
$ws.Range("I97").Value = 1027.5588
$ws.Range("K97").Value = 1027.5588
$ws.Range("M97").Value = -531.5588
$ws.Range("H132").Value = 4669.264
$ws.Range("I132").Value = 4360.2705
$ws.Range("J132").Value = 5383.8125
$ws.Range("K132").Value = 13080.8115
$ws.Range("L132").Value = 16151.4375
$ws.Range("M132").Value = -10550.8115
$ws.Range("N132").Value = -21211.4375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 13203428
$ws.Range("I22").Value = 16502422
$ws.Range("J22").Value = 7450
$ws.Range("K22").Value = 16502422
$ws.Range("L22").Value = 7450
$ws.Range("M22").Value = -16502127
$ws.Range("N22").Value = -8040
$ws.Range("H27").Value = 13203428
$ws.Range("I27").Value = 16502422
$ws.Range("J27").Value = 7450
$ws.Range("K27").Value = 16502422
$ws.Range("L27").Value = 7450
$ws.Range("M27").Value = -16502315
$ws.Range("N27").Value = -7664
$ws.Range("H46").Value = 1333.3334
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -2376
$ws.Range("H93").Value = 2647254.5
$ws.Range("I93").Value = 624.63635
$ws.Range("K93").Value = 624.63635
$ws.Range("M93").Value = 623.36365
$ws.Range("H132").Value = 3441.9268
$ws.Range("I132").Value = 2624.2812
$ws.Range("J132").Value = 6349.1113
$ws.Range("K132").Value = 7872.8436
$ws.Range("L132").Value = 19047.3339
$ws.Range("M132").Value = -5342.8436
$ws.Range("N132").Value = -24107.3339
$ws.Range("H136").Value = 4543
$ws.Range("I136").Value = 4000.4
$ws.Range("K136").Value = 12001.2
$ws.Range("M136").Value = -9451.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 63333
$ws.Range("J130").Value = 63333
$ws.Range("L130").Value = 63333
$ws.Range("N130").Value = -73373
$ws.Range("H132").Value = 1355.9048
$ws.Range("I132").Value = 1048.725
$ws.Range("K132").Value = 3146.175
$ws.Range("M132").Value = -616.1749999999997
